$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103, shifting all existing rows (103..205) down by one.
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with the new weekly record.
$ws.Cells.Item(103, 1).Value = 3
$ws.Cells.Item(103, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(103, 3).Value = "Coquimbo"
$ws.Cells.Item(103, 4).Value = 44629
$ws.Cells.Item(103, 5).Value = 5
$ws.Cells.Item(103, 6).Value = 100112010
$ws.Cells.Item(103, 7).Value = "Achicoria"
$ws.Cells.Item(103, 8).Value = "Sin especificar"
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value = 60
$ws.Cells.Item(103, 11).Value = 7500
$ws.Cells.Item(103, 12).Value = 7500
$ws.Cells.Item(103, 13).Value = 7500
$ws.Cells.Item(103, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(103, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(103, 16).Value = 469
$ws.Cells.Item(103, 17).Value = 16
$ws.Cells.Item(103, 18).Value = "Hortaliza"
